$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B5").Value = 0.05
$wsSummary.Range("B6").Value = 11
$wsSummary.Range("B9").Value = 36.36

# --- Strategy Status sheet ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("D4").Value = 11
$wsStatus.Range("G4").Value = 36.36

# --- New trade row data (Trade #11) ---
$newRow = @{
    A = 11
    B = "2026-02-17"
    C = "04:07:05"
    D = "MarketMaking"
    E = "UP"
    F = 0.21
    G = 0.21
    H = "CLOSED"
    I = 0
    J = 0
    K = 100.03
    L = 0
    M = 0
    N = 0.6
    O = "Normal spread capture: 19600 bps"
    P = "early_exit"
    Q = 0.12
}

function Add-TradeRow($ws, $rowNum, $data) {
    $ws.Cells.Item($rowNum, 1).Value = $data.A
    $ws.Cells.Item($rowNum, 2).NumberFormat = "@"
    $ws.Cells.Item($rowNum, 2).Value = $data.B
    $ws.Cells.Item($rowNum, 3).NumberFormat = "@"
    $ws.Cells.Item($rowNum, 3).Value = $data.C
    $ws.Cells.Item($rowNum, 4).Value = $data.D
    $ws.Cells.Item($rowNum, 5).Value = $data.E
    $ws.Cells.Item($rowNum, 6).Value = $data.F
    $ws.Cells.Item($rowNum, 7).Value = $data.G
    $ws.Cells.Item($rowNum, 8).Value = $data.H
    $ws.Cells.Item($rowNum, 9).Value = $data.I
    $ws.Cells.Item($rowNum, 10).Value = $data.J
    $ws.Cells.Item($rowNum, 11).Value = $data.K
    $ws.Cells.Item($rowNum, 12).Value = $data.L
    $ws.Cells.Item($rowNum, 13).Value = $data.M
    $ws.Cells.Item($rowNum, 14).Value = $data.N
    $ws.Cells.Item($rowNum, 15).Value = $data.O
    $ws.Cells.Item($rowNum, 16).Value = $data.P
    $ws.Cells.Item($rowNum, 17).Value = $data.Q
}

# --- All Trades sheet ---
$wsAllTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $wsAllTrades 12 $newRow

# --- MarketMaking sheet ---
$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $wsMarketMaking 12 $newRow
